$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.065754003257296
$ws.Range("D2").Value = 1.075513797482597
$ws.Range("E2").Value = 1.071381818570679
$ws.Range("F2").Value = 1.083315849450454
$ws.Range("I2").Value = 1.049096248622849
$ws.Range("J2").Value = 1.070707502821462
$ws.Range("K2").Value = 1.078200167420375
$ws.Range("L2").Value = 1.074079141773484
$ws.Range("M2").Value = 1.08598178331625
$ws.Range("N2").Value = 1.02717197696528
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.067039899323876
$ws.Range("D3").Value = 1.07672745977088
$ws.Range("E3").Value = 1.072527167338824
$ws.Range("F3").Value = 1.084549191914757
$ws.Range("I3").Value = 1.049398172960842
$ws.Range("J3").Value = 1.071647532168435
$ws.Range("K3").Value = 1.079230188032737
$ws.Range("L3").Value = 1.075040221041172
$ws.Range("M3").Value = 1.087032921357636
$ws.Range("N3").Value = 1.02749133821702
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.067871415466116
$ws.Range("D4").Value = 1.077512511587654
$ws.Range("E4").Value = 1.073268086722166
$ws.Range("F4").Value = 1.085347038162497
$ws.Range("I4").Value = 1.049591457502019
$ws.Range("J4").Value = 1.072254723255136
$ws.Range("K4").Value = 1.079895837857626
$ws.Range("L4").Value = 1.075661319218596
$ws.Range("M4").Value = 1.087712297195847
$ws.Range("N4").Value = 1.027697469756699
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.06822085724716
$ws.Range("D5").Value = 1.07784248503117
$ws.Range("E5").Value = 1.073579523650474
$ws.Range("F5").Value = 1.085682404738006
$ws.Range("I5").Value = 1.049672216772912
$ws.Range("J5").Value = 1.07250973178019
$ws.Range("K5").Value = 1.080175477200263
$ws.Range("L5").Value = 1.075922242778318
$ws.Range("M5").Value = 1.087997721685946
$ws.Range("N5").Value = 1.027784004303999
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.068279522657399
$ws.Range("D6").Value = 1.077897885427157
$ws.Range("E6").Value = 1.073631812660028
$ws.Range("F6").Value = 1.08573871150566
$ws.Range("I6").Value = 1.049685747449751
$ws.Range("J6").Value = 1.072552533953525
$ws.Range("K6").Value = 1.080222418210293
$ws.Range("L6").Value = 1.075966042126569
$ws.Range("M6").Value = 1.08804563491092
$ws.Range("N6").Value = 1.027798526628338
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.067876085222253
$ws.Range("D7").Value = 1.077516920951408
$ws.Range("E7").Value = 1.073272248334531
$ws.Range("F7").Value = 1.085351519533172
$ws.Range("I7").Value = 1.049592538565686
$ws.Range("J7").Value = 1.072258131689822
$ws.Range("K7").Value = 1.079899575195202
$ws.Range("L7").Value = 1.075664806421213
$ws.Range("M7").Value = 1.087716111776989
$ws.Range("N7").Value = 1.027698626518828
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.066188692288237
$ws.Range("D8").Value = 1.075924017219614
$ws.Range("E8").Value = 1.071768936595649
$ws.Range("F8").Value = 1.083732707750741
$ws.Range("I8").Value = 1.049198716535099
$ws.Range("J8").Value = 1.071025412895092
$ws.Range("K8").Value = 1.078548443232625
$ws.Range("L8").Value = 1.074404106427835
$ws.Range("M8").Value = 1.086337183089768
$ws.Range("N8").Value = 1.027280013857643
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.063210994615946
$ws.Range("D9").Value = 1.073114967097621
$ws.Range("E9").Value = 1.069118313132635
$ws.Range("F9").Value = 1.080878465088139
$ws.Range("I9").Value = 1.0484887899057
$ws.Range("J9").Value = 1.068844925771135
$ws.Range("K9").Value = 1.076161041920498
$ws.Range("L9").Value = 1.072176509269538
$ws.Range("M9").Value = 1.083901273314458
$ws.Range("N9").Value = 1.026538385644028
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.061222797548233
$ws.Range("D10").Value = 1.071240690722992
$ws.Range("E10").Value = 1.067350041851489
$ws.Range("F10").Value = 1.07897437134853
$ws.Range("I10").Value = 1.048004734958331
$ws.Range("J10").Value = 1.067385589007312
$ws.Range("K10").Value = 1.07456492968041
$ws.Range("L10").Value = 1.070687246875739
$ws.Range("M10").Value = 1.08227314130973
$ws.Range("N10").Value = 1.026041258109073
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.060361119003976
$ws.Range("D11").Value = 1.070428704491926
$ws.Range("E11").Value = 1.066584052179851
$ws.Range("F11").Value = 1.078149547882907
$ws.Range("I11").Value = 1.047792569985361
$ws.Range("J11").Value = 1.066752307351756
$ws.Range("K11").Value = 1.073872698256422
$ws.Range("L11").Value = 1.070041359356678
$ws.Range("M11").Value = 1.081567120942013
$ws.Range("N11").Value = 1.025825346188435
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.060040933057453
$ws.Range("D12").Value = 1.070127032031049
$ws.Range("E12").Value = 1.066299479550659
$ws.Range("F12").Value = 1.077843118221502
$ws.Range("I12").Value = 1.047713376077178
$ws.Range("J12").Value = 1.066516868807483
$ws.Range("K12").Value = 1.073615404327781
$ws.Range("L12").Value = 1.069801291425873
$ws.Range("M12").Value = 1.081304716598055
$ws.Range("N12").Value = 1.025745048160934
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.060109619525283
$ws.Range("D13").Value = 1.070191744765628
$ws.Range("E13").Value = 1.066360523645149
$ws.Range("F13").Value = 1.077908850888218
$ws.Range("I13").Value = 1.047730380938532
$ws.Range("J13").Value = 1.066567380694173
$ws.Range("K13").Value = 1.073670602406839
$ws.Range("L13").Value = 1.069852793932237
$ws.Range("M13").Value = 1.081361010361789
$ws.Range("N13").Value = 1.025762276839426
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.060334654805872
$ws.Range("D14").Value = 1.070403769486091
$ws.Range("E14").Value = 1.066560530350661
$ws.Range("F14").Value = 1.078124219418433
$ws.Range("I14").Value = 1.047786031679581
$ws.Range("J14").Value = 1.066732850216336
$ws.Range("K14").Value = 1.073851433708652
$ws.Range("L14").Value = 1.070021518473177
$ws.Range("M14").Value = 1.081545433732793
$ws.Range("N14").Value = 1.025818710746541
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.060473290438711
$ws.Range("D15").Value = 1.070534396270292
$ws.Range("E15").Value = 1.066683754390715
$ws.Range("F15").Value = 1.078256907908286
$ws.Range("I15").Value = 1.047820268710857
$ws.Range("J15").Value = 1.066834773614002
$ws.Range("K15").Value = 1.073962827459533
$ws.Range("L15").Value = 1.070125454420021
$ws.Range("M15").Value = 1.081659042188909
$ws.Range("N15").Value = 1.025853468435753
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.061279967659655
$ws.Range("D16").Value = 1.071294570660977
$ws.Range("E16").Value = 1.067400871189404
$ws.Range("F16").Value = 1.079029104809808
$ws.Range("I16").Value = 1.048018761511545
$ws.Range("J16").Value = 1.067427588568267
$ws.Range("K16").Value = 1.074610847325395
$ws.Range("L16").Value = 1.070730090472896
$ws.Range("M16").Value = 1.082319975680802
$ws.Range("N16").Value = 1.026055573666116
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.061785764670734
$ws.Range("D17").Value = 1.071771295720111
$ws.Range("E17").Value = 1.067850613172735
$ws.Range("F17").Value = 1.079513391066712
$ws.Range("I17").Value = 1.048142583195244
$ws.Range("J17").Value = 1.067799074833423
$ws.Range("K17").Value = 1.075017035692229
$ws.Range("L17").Value = 1.071109086138881
$ws.Range("M17").Value = 1.082734284739193
$ws.Range("N17").Value = 1.026182173768103
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.062080712871275
$ws.Range("D18").Value = 1.072049321795017
$ws.Range("E18").Value = 1.068112909685564
$ws.Range("F18").Value = 1.079795834587491
$ws.Range("I18").Value = 1.048214558768517
$ws.Range("J18").Value = 1.068015623572629
$ws.Range("K18").Value = 1.075253852113688
$ws.Range("L18").Value = 1.071330048954122
$ws.Range("M18").Value = 1.082975845137384
$ws.Range("N18").Value = 1.026255954565348
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.062181270032766
$ws.Range("D19").Value = 1.072144114867909
$ws.Range("E19").Value = 1.068202340937624
$ws.Range("F19").Value = 1.079892135103535
$ws.Range("I19").Value = 1.048239058621562
$ws.Range("J19").Value = 1.06808943862188
$ws.Range("K19").Value = 1.075334582422909
$ws.Range("L19").Value = 1.071405374848391
$ws.Range("M19").Value = 1.083058194260151
$ws.Range("N19").Value = 1.026281101262758
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.06173150518565
$ws.Range("D20").Value = 1.071720151748251
$ws.Range("E20").Value = 1.067802363274776
$ws.Range("F20").Value = 1.079461435112559
$ws.Range("I20").Value = 1.04812932390942
$ws.Range("J20").Value = 1.067759231638704
$ws.Range("K20").Value = 1.074973466546856
$ws.Range("L20").Value = 1.071068433742234
$ws.Range("M20").Value = 1.082689843563734
$ws.Range("N20").Value = 1.026168597286846
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.060268390870604
$ws.Range("D21").Value = 1.070341335296066
$ws.Range("E21").Value = 1.066501634759428
$ws.Range("F21").Value = 1.078060800252651
$ws.Range("I21").Value = 1.047769654596786
$ws.Range("J21").Value = 1.066684129369413
$ws.Range("K21").Value = 1.07379818806457
$ws.Range("L21").Value = 1.069971837640787
$ws.Range("M21").Value = 1.081491129998665
$ws.Range("N21").Value = 1.025802095098526
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.059347775789179
$ws.Range("D22").Value = 1.069474043444475
$ws.Range("E22").Value = 1.065683524864903
$ws.Range("F22").Value = 1.077179854648719
$ws.Range("I22").Value = 1.047541280245423
$ws.Range("J22").Value = 1.066006956400384
$ws.Range("K22").Value = 1.073058268555292
$ws.Range("L22").Value = 1.069281457288841
$ws.Range("M22").Value = 1.080736541960926
$ws.Range("N22").Value = 1.025571089069959
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.059835878529414
$ws.Range("D23").Value = 1.069933847734497
$ws.Range("E23").Value = 1.066117248717262
$ws.Range("F23").Value = 1.07764689094518
$ws.Range("I23").Value = 1.047662558051031
$ws.Range("J23").Value = 1.066366054411028
$ws.Range("K23").Value = 1.073450607153261
$ws.Range("L23").Value = 1.069647527687099
$ws.Range("M23").Value = 1.08113665035694
$ws.Range("N23").Value = 1.025693604150897
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.061756022941911
$ws.Range("D24").Value = 1.071743261629968
$ws.Range("E24").Value = 1.067824165420632
$ws.Range("F24").Value = 1.079484911872108
$ws.Range("I24").Value = 1.048135315974748
$ws.Range("J24").Value = 1.067777235474022
$ws.Range("K24").Value = 1.074993153896913
$ws.Range("L24").Value = 1.071086803116123
$ws.Range("M24").Value = 1.082709924924057
$ws.Range("N24").Value = 1.02617473210869
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.063981328082024
$ws.Range("D25").Value = 1.073841442090598
$ws.Range("E25").Value = 1.069803763855253
$ws.Range("F25").Value = 1.081616568709274
$ws.Range("I25").Value = 1.048674218467169
$ws.Range("J25").Value = 1.069409626584503
$ws.Range("K25").Value = 1.076779028830543
$ws.Range("L25").Value = 1.072753128428622
$ws.Range("M25").Value = 1.084531744652643
$ws.Range("N25").Value = 1.026730589287099
